$d = $word.ActiveDocument

# --- 1. The "m: link 'bookmark1' 'a reference to bookmark1'" field (a broken/
#        unresolved "link" M2Doc query) is rewritten by the new
#        TokenIteratorFieldRewriterSplit parser as plain literal text, i.e. the
#        field code runs (fldChar begin / instrText* / fldChar end) disappear
#        and are replaced by the same text, now wrapped in "{" / "}", stored as
#        normal w:t runs. ---
$targetField = $null
for ($i = 1; $i -le $d.Fields.Count; $i++) {
    $candidate = $d.Fields.Item($i)
    if ($candidate.Code.Text -eq "m: link 'bookmark1' 'a reference to bookmark1'") {
        $targetField = $candidate
    }
}

if ($targetField -ne $null) {
    $insertPos = $targetField.Code.Start - 1
    $targetField.Delete()

    $gap = $d.Range($insertPos, $insertPos)
    $gap.InsertAfter("{m: link 'bookmark1' 'a reference to bookmark1'}")
}

# --- 2. Flag the broken-link error message so it stands out after the now
#        literal "{m: ...}" text: "Couldn't find the 'link' variable" becomes
#        "    <---Couldn't find the 'link' variable". ---
$errRange = $d.Content
$errRange.Find.Execute("Couldn't find the 'link' variable") | Out-Null
if ($errRange.Find.Found) {
    $errRange.Text = "    <---Couldn't find the 'link' variable"
}
